# Implement database schema migration: append a new row (row 37) of
# sensor-log data to each of the four worksheets, mirroring the existing
# row-36 layout/formatting.

$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{
        Name = "FE_LFT_#1"
        A = 45823.4937037037
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x68"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 360
        I = 15
    },
    @{
        Name = "FE_LFT_#2"
        A = 45823.4937037037
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x7C"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 380
        I = 14
    },
    @{
        Name = "FE_PLT_#1"
        A = 45823.4937037037
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x6B"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 107
        I = 3
    },
    @{
        Name = "FE_PLT_#2"
        A = 45823.4937037037
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x6B"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 107
        I = 3
    }
)

foreach ($sheetData in $sheetsData) {
    $ws = $wb.Worksheets.Item($sheetData.Name)
    $newRow = 37

    # New A cell: copy the date/time number format from the row above (A36)
    $ws.Range("A$newRow").Value = $sheetData.A
    $ws.Range("A$newRow").NumberFormat = $ws.Range("A36").NumberFormat

    $ws.Range("B$newRow").Value = $sheetData.B
    $ws.Range("C$newRow").Value = $sheetData.C
    $ws.Range("D$newRow").Value = $sheetData.D
    $ws.Range("E$newRow").Value = $sheetData.E
    $ws.Range("F$newRow").Value = $sheetData.F
    $ws.Range("G$newRow").Value = $sheetData.G
    $ws.Range("H$newRow").Value = $sheetData.H
    $ws.Range("I$newRow").Value = $sheetData.I
}
